$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Switch the standard-deviation formula in C6 from the modern STDEV.P
# (population stdev) to the legacy STDEV (sample stdev) function.
$ws.Range("C6").Formula = "=STDEV('Konzentration des Blei Messgerä'!G18:G20)"

# Move the active selection from E12 to C7.
$ws.Range("C7").Select()
